$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to Text
# so Excel stores them as strings (matching the source data), not as numbers.
$textCells = @(
    "D5", "D6", "D7", "D10", "D11", "D12", "D16", "D18",
    "D19", "D22", "D23", "D24", "D25", "D27", "D28", "D29",
    "D30", "D31", "D32", "D33", "D34", "D35", "D38", "D42",
    "D43", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '46.817.06'
$ws.Range("E2").Value = '  +4.07%  '
$ws.Range("D3").Value = '2.270.14'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '302.55'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '101.39'
$ws.Range("E6").Value = '  +7.08%  '
$ws.Range("D7").Value = '0.563'
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").Value = '35.92'
$ws.Range("E10").Value = '  +4.82%  '
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '7.19'
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '2.618.85'
$ws.Range("E14").Value = '  +0.39%  '
$ws.Range("D15").Value = '2.274.45'
$ws.Range("E15").Value = '  +0.51%  '
$ws.Range("D16").Value = '13.67'
$ws.Range("E16").Value = '  +0.33%  '
$ws.Range("D17").Value = '46.787.04'
$ws.Range("E17").Value = '  +4.32%  '
$ws.Range("D18").Value = '0.799'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = '13.09'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").Value = '0.0₃0932'
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("E21").Value = '  -2.57%  '
$ws.Range("D22").Value = '65.47'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '250.56'
$ws.Range("E23").Value = '  +4.92%  '
$ws.Range("D24").Value = '2.85'
$ws.Range("E24").Value = '  -1.74%  '
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = '43.07'
$ws.Range("E27").Value = '  +4.01%  '
$ws.Range("D28").Value = '2.25'
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("D29").Value = '9.76'
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").Value = '19.88'
$ws.Range("E30").Value = '  +1.67%  '
$ws.Range("D31").Value = '2.79'
$ws.Range("E31").Value = '  +8.49%  '
$ws.Range("B32").Value = 'Monero'
$ws.Range("C32").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D32").Value = '146.98'
$ws.Range("E32").Value = '  -4.02%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '5.48'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").Value = '3.26'
$ws.Range("E34").Value = '  +9.97%  '
$ws.Range("D35").Value = '0.0776'
$ws.Range("E35").Value = '  -1.38%  '
$ws.Range("E36").Value = '  +10.36%  '
$ws.Range("E37").Value = '  -0.91%  '
$ws.Range("D38").Value = '16.47'
$ws.Range("E38").Value = '  +20.48%  '
$ws.Range("E39").Value = '  -1.80%  '
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("E41").Value = '  -3.90%  '
$ws.Range("D42").Value = '3.23'
$ws.Range("E42").Value = '  -0.56%  '
$ws.Range("D43").Value = '0.999'
$ws.Range("E44").Value = '  +2.35%  '
$ws.Range("D45").Value = '1.808.49'
$ws.Range("E45").Value = '  +3.18%  '
$ws.Range("D46").Value = '91.41'
$ws.Range("E46").Value = '  +21.26%  '
$ws.Range("D47").Value = '0.190'
$ws.Range("E47").Value = '  -3.57%  '
$ws.Range("D48").Value = '72.61'
$ws.Range("E48").Value = '  +2.38%  '
$ws.Range("D49").Value = '4.83'
$ws.Range("E49").Value = '  +3.69%  '
$ws.Range("D50").Value = '94.24'
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").Value = '7.84'
$ws.Range("E51").Value = '  +0.47%  '
